$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D17").Value = "0.14027746021747589"
$ws.Range("E17").Value = "0.11668840050697327"
$ws.Range("H17").Value = "-0.088431805372238159"
$ws.Range("I17").Value = "0.36898672580718994"

$ws.Range("D18").Value = "0.10525839775800705"
$ws.Range("E18").Value = "0.11187935620546341"
$ws.Range("H18").Value = "-0.11402513831853867"
$ws.Range("I18").Value = "0.32454192638397217"

$ws.Range("D19").Value = "0.086012274026870728"
$ws.Range("E19").Value = "0.10779380798339844"
$ws.Range("H19").Value = "-0.12526358664035797"
$ws.Range("I19").Value = "0.29728814959526062"

$ws.Range("D20").Value = "0.77207005023956299"
$ws.Range("E20").Value = "0.042151398956775665"
$ws.Range("H20").Value = "0.68945330381393433"
$ws.Range("I20").Value = "0.85468679666519165"

$ws.Range("D21").Value = "0.76392996311187744"
$ws.Range("E21").Value = "0.042095102369785309"
$ws.Range("H21").Value = "0.68142354488372803"
$ws.Range("I21").Value = "0.84643638134002686"

$ws.Range("D22").Value = "0.7738615870475769"
$ws.Range("E22").Value = "0.0384235680103302"
$ws.Range("H22").Value = "0.69855141639709473"
$ws.Range("I22").Value = "0.84917175769805908"

$ws.Range("D23").Value = "1.0425946712493896"
$ws.Range("E23").Value = "0.11899025738239288"
$ws.Range("H23").Value = "0.80937379598617554"
$ws.Range("I23").Value = "1.2758156061172485"

$ws.Range("D24").Value = "1.0453108549118042"
$ws.Range("E24").Value = "0.11903690546751022"
$ws.Range("H24").Value = "0.81199854612350464"
$ws.Range("I24").Value = "1.2786232233047485"

$ws.Range("D25").Value = "1.0401633977890015"
$ws.Range("E25").Value = "0.11633656919002533"
$ws.Range("H25").Value = "0.81214374303817749"
$ws.Range("I25").Value = "1.2681831121444702"

$ws.Range("D26").Value = "0.13906879723072052"
$ws.Range("E26").Value = "0.094747200608253479"
$ws.Range("H26").Value = "-0.046635717153549194"
$ws.Range("I26").Value = "0.32477331161499023"

$ws.Range("D27").Value = "0.10438137501478195"
$ws.Range("E27").Value = "0.088896416127681732"
$ws.Range("H27").Value = "-0.069855600595474243"
$ws.Range("I27").Value = "0.27861833572387695"

$ws.Range("D28").Value = "0.10824931412935257"
$ws.Range("E28").Value = "0.088015884160995483"
$ws.Range("H28").Value = "-0.064261816442012787"
$ws.Range("I28").Value = "0.28076043725013733"

$ws.Range("D29").Value = "0.26149716973304749"
$ws.Range("E29").Value = "0.17170630395412445"
$ws.Range("H29").Value = "-0.075047187507152557"
$ws.Range("I29").Value = "0.59804153442382812"

$ws.Range("D30").Value = "0.22199662029743195"
$ws.Range("E30").Value = "0.17070366442203522"
$ws.Range("H30").Value = "-0.11258256435394287"
$ws.Range("I30").Value = "0.55657577514648438"

$ws.Range("D31").Value = "0.19604043662548065"
$ws.Range("E31").Value = "0.16363796591758728"
$ws.Range("H31").Value = "-0.1246899738907814"
$ws.Range("I31").Value = "0.51677083969116211"

$ws.Range("D32").Value = "0.004408013541251421"
$ws.Range("E32").Value = "0.0013329912908375263"
$ws.Range("H32").Value = "0.0017953505739569664"
$ws.Range("I32").Value = "0.0070206765085458755"

$ws.Range("D33").Value = "0.0039839497767388821"
$ws.Range("E33").Value = "0.0016661899862810969"
$ws.Range("H33").Value = "0.00071821740129962564"
$ws.Range("I33").Value = "0.0072496822103857994"

$ws.Range("D34").Value = "0.0039056963287293911"
$ws.Range("E34").Value = "0.0015684141544625163"
$ws.Range("H34").Value = "0.00083160458598285913"
$ws.Range("I34").Value = "0.0069797881878912449"

$ws.Range("D35").Value = "0.25779908895492554"
$ws.Range("E35").Value = "0.041645646095275879"
$ws.Range("H35").Value = "0.1761736273765564"
$ws.Range("I35").Value = "0.33942455053329468"

$ws.Range("D36").Value = "0.25669911503791809"
$ws.Range("E36").Value = "0.042932983487844467"
$ws.Range("H36").Value = "0.17255046963691711"
$ws.Range("I36").Value = "0.34084776043891907"

$ws.Range("D37").Value = "0.25736364722251892"
$ws.Range("E37").Value = "0.042606338858604431"
$ws.Range("H37").Value = "0.17385523021221161"
$ws.Range("I37").Value = "0.34087207913398743"

$ws.Range("D38").Value = "5.0924191474914551"
$ws.Range("E38").Value = "0.57012712955474854"
$ws.Range("H38").Value = "3.9749698638916016"
$ws.Range("I38").Value = "6.2098684310913086"

$ws.Range("D39").Value = "5.1330761909484863"
$ws.Range("E39").Value = "0.5950850248336792"
$ws.Range("H39").Value = "3.9667096138000488"
$ws.Range("I39").Value = "6.2994427680969238"

$ws.Range("D40").Value = "5.1145944595336914"
$ws.Range("E40").Value = "0.57889062166213989"
$ws.Range("H40").Value = "3.979968786239624"
$ws.Range("I40").Value = "6.2492198944091797"

$ws.Range("D41").Value = "0.0014742759522050619"
$ws.Range("E41").Value = "0.00070770381717011333"
$ws.Range("H41").Value = "0.000087176471424754709"
$ws.Range("I41").Value = "0.002861375454813242"

$ws.Range("D42").Value = "0.0012385620502755046"
$ws.Range("E42").Value = "0.0006318475934676826"
$ws.Range("H42").Value = "0.00000014076708509946911"
$ws.Range("I42").Value = "0.0024769832380115986"

$ws.Range("D43").Value = "0.0012464105384424329"
$ws.Range("E43").Value = "0.00060741853667423129"
$ws.Range("H43").Value = "0.000055870204960228875"
$ws.Range("I43").Value = "0.0024369508028030396"

$ws.Range("D44").Value = "0.013377105817198753"
$ws.Range("E44").Value = "0.0069817281328141689"
$ws.Range("H44").Value = "-0.00030708132544532418"
$ws.Range("I44").Value = "0.02706129290163517"

$ws.Range("D45").Value = "0.011803300119936466"
$ws.Range("E45").Value = "0.0081560183316469193"
$ws.Range("H45").Value = "-0.0041824956424534321"
$ws.Range("I45").Value = "0.027789095416665077"

$ws.Range("D46").Value = "0.010615397244691849"
$ws.Range("E46").Value = "0.00759173184633255"
$ws.Range("H46").Value = "-0.0042643970809876919"
$ws.Range("I46").Value = "0.025495192036032677"

